{"js": "const body = context.document.body;\n\n// 1) \"2. Analiza SWOT (Mocne i s\u0142abe strony projektu)\" -> \"2. Analiza SWOT \"\nconst swotResults = body.search(\"2. Analiza SWOT (Mocne i s\u0142abe strony projektu)\", { matchCase: true });\nswotResults.load(\"items\");\nawait context.sync();\nfor (const r of swotResults.items) {\n  r.insertText(\"2. Analiza SWOT \", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"Zmiana podej\u015bcia do magazynu (Najwi\u0119ksza zmiana)\" -> \"Zmiana podej\u015bcia do magazynu\"\nconst magResults = body.search(\"Zmiana podej\u015bcia do magazynu (Najwi\u0119ksza zmiana)\", { matchCase: true });\nmagResults.load(\"items\");\nawait context.sync();\nfor (const r of magResults.items) {\n  r.insertText(\"Zmiana podej\u015bcia do magazynu\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Remove the unnecessary trailing parentheticals from two headings:\n#   \"2. Analiza SWOT (Mocne i s\u0142abe strony projektu)\" -> \"2. Analiza SWOT \"\n#   \"Zmiana podej\u015bcia do magazynu (Najwi\u0119ksza zmiana)\" -> \"Zmiana podej\u015bcia do magazynu\"\n#\n# In both cases the replacement text is a literal prefix of the original text,\n# so we overwrite the Range with FormattedText copied from the leading\n# sub-range we want to keep (rather than assigning .Text directly), which\n# preserves the run's existing formatting/identity instead of fabricating a\n# brand-new run.\n\n$d = $word.ActiveDocument\n\nfunction Remove-TrailingParenthetical($doc, [string]$oldText, [string]$newText) {\n    $content = $doc.Content.Text\n    $start = $content.IndexOf($oldText)\n    if ($start -lt 0) {\n        throw \"Text not found: $oldText\"\n    }\n    $target = $doc.Range($start, $start + $oldText.Length)\n    $keep = $doc.Range($start, $start + $newText.Length)\n    $target.FormattedText = $keep.FormattedText\n}\n\nRemove-TrailingParenthetical $d \"2. Analiza SWOT (Mocne i s\u0142abe strony projektu)\" \"2. Analiza SWOT \"\nRemove-TrailingParenthetical $d \"Zmiana podej\u015bcia do magazynu (Najwi\u0119ksza zmiana)\" \"Zmiana podej\u015bcia do magazynu\"\n"}
